$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# need an explicit Text number format so the literal string (incl. trailing
# zeros / decimal formatting) is preserved, matching the source data exactly.

$ws.Range("D2").Value = "70.381.52"
$ws.Range("E2").Value = "  +0.76%  "

$ws.Range("D3").Value = "3.757.72"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.59"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.76"
$ws.Range("E6").Value = "  +1.96%  "

$ws.Range("D7").Value = "3.753.00"
$ws.Range("E7").Value = "  -0.83%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.38"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("E12").Value = "  -3.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.19"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000256"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").Value = "4.388.06"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("D16").Value = "3.755.46"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").Value = "70.615.64"
$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.121"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.52"
$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "506.61"
$ws.Range("E21").Value = "  -2.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.24"
$ws.Range("E22").Value = "  -1.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  -2.52%  "

$ws.Range("E24").Value = "  +4.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.37"
$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.08"
$ws.Range("E26").Value = "  -3.19%  "

$ws.Range("E27").Value = "  +3.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000135"
$ws.Range("E28").Value = "  +9.23%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("E30").Value = "  -1.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.93"
$ws.Range("E31").Value = "  +2.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.87"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.59"
$ws.Range("E33").Value = "  -4.21%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.11"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.351"
$ws.Range("E38").Value = "  +2.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  +5.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.14"
$ws.Range("E40").Value = "  +17.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.08"
$ws.Range("E41").Value = "  -4.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "50.08"
$ws.Range("E42").Value = "  -2.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "437.34"
$ws.Range("E43").Value = "  +2.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.36"
$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.63"
$ws.Range("E45").Value = "  -2.28%  "

$ws.Range("D46").Value = "2.962.03"
$ws.Range("E46").Value = "  -4.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0363"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.22"
$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.70"
$ws.Range("E50").Value = "  -0.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").Value = "  -0.82%  "
